$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows: set Price (D) and Volume(1h) (E) values, and for row 28-51
# also update Coin (B) and Link (C) due to the new "LEO" entry shifting the list down ---

# Cells whose new value is a plain parseable number but must remain TEXT
# (matching the original inline-string formatting of the Price column).
# We temporarily force a text number-format, assign the values, then restore
# the format back to General so no stray formatting is left behind.
$forceTextRange = $ws.Range('D4,D5,D6,D10,D11,D12,D14,D15,D17,D19,D20,D22,D23,D24,D27,D28,D29,D30,D31,D32,D33,D34,D35,D36,D37,D38,D39,D40,D41,D42,D43,D45,D46,D47,D48,D49,D50')
$forceTextRange.NumberFormat = "@"

$ws.Range('D4').Value = '0.999'
$ws.Range('D5').Value = '306.66'
$ws.Range('D6').Value = '101.11'
$ws.Range('D10').Value = '34.85'
$ws.Range('D11').Value = '52.49'
$ws.Range('D12').Value = '0.0798'
$ws.Range('D14').Value = '6.87'
$ws.Range('D15').Value = '15.85'
$ws.Range('D17').Value = '0.818'
$ws.Range('D19').Value = '6.24'
$ws.Range('D20').Value = '11.81'
$ws.Range('D22').Value = '67.94'
$ws.Range('D23').Value = '237.02'
$ws.Range('D24').Value = '2.03'
$ws.Range('D27').Value = '25.45'
$ws.Range('D28').Value = '3.95'
$ws.Range('D29').Value = '2.32'
$ws.Range('D30').Value = '35.36'
$ws.Range('D31').Value = '9.41'
$ws.Range('D32').Value = '164.03'
$ws.Range('D33').Value = '0.999'
$ws.Range('D34').Value = '5.12'
$ws.Range('D35').Value = '17.56'
$ws.Range('D36').Value = '4.62'
$ws.Range('D37').Value = '0.0728'
$ws.Range('D38').Value = '2.45'
$ws.Range('D39').Value = '1.86'
$ws.Range('D40').Value = '2.92'
$ws.Range('D41').Value = '0.102'
$ws.Range('D42').Value = '0.113'
$ws.Range('D43').Value = '2.57'
$ws.Range('D45').Value = '0.0286'
$ws.Range('D46').Value = '18.82'
$ws.Range('D47').Value = '10.17'
$ws.Range('D48').Value = '2.93'
$ws.Range('D49').Value = '56.39'
$ws.Range('D50').Value = '2.89'

$forceTextRange.NumberFormat = "General"

# Remaining cells (strings that Excel will not auto-convert to numbers)
$ws.Range('D2').Value = '42.986.52'
$ws.Range('E2').Value = '  -0.97%  '
$ws.Range('D3').Value = '2.337.14'
$ws.Range('E3').Value = '  +1.23%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('E5').Value = '  -1.35%  '
$ws.Range('E6').Value = '  -2.06%  '
$ws.Range('E7').Value = '  -4.03%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -3.76%  '
$ws.Range('E10').Value = '  -2.50%  '
$ws.Range('E11').Value = '  +1.02%  '
$ws.Range('E12').Value = '  -1.95%  '
$ws.Range('E13').Value = '  +0.74%  '
$ws.Range('E14').Value = '  -2.22%  '
$ws.Range('E15').Value = '  +5.46%  '
$ws.Range('D16').Value = '2.351.98'
$ws.Range('E16').Value = '  -2.55%  '
$ws.Range('E17').Value = '  +1.17%  '
$ws.Range('D18').Value = '42.913.65'
$ws.Range('E18').Value = '  -0.90%  '
$ws.Range('E19').Value = '  +0.88%  '
$ws.Range('E20').Value = '  -3.88%  '
$ws.Range('E21').Value = '  -2.60%  '
$ws.Range('E22').Value = '  -0.30%  '
$ws.Range('E23').Value = '  -1.91%  '
$ws.Range('E24').Value = '  +0.89%  '
$ws.Range('E25').Value = '  -2.07%  '
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('E27').Value = '  +1.78%  '
$ws.Range('B28').Value = 'LEO'
$ws.Range('C28').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('E28').Value = '  -0.77%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('E29').Value = '  +1.07%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('E30').Value = '  -4.02%  '
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('E31').Value = '  -2.70%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('E32').Value = '  -4.55%  '
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('E34').Value = '  -3.06%  '
$ws.Range('B35').Value = 'Celestia'
$ws.Range('C35').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('E35').Value = '  -1.67%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('E36').Value = '  +6.52%  '
$ws.Range('E37').Value = '  -1.90%  '
$ws.Range('B38').Value = 'WEMIXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('E38').Value = '  -4.05%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('E39').Value = '  -1.51%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('E40').Value = '  -4.94%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('E41').Value = '  -4.80%  '
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('E42').Value = '  -2.47%  '
$ws.Range('B43').Value = 'ApeXProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('E43').Value = '  +11.38%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.023.80'
$ws.Range('E44').Value = '  +2.54%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('E45').Value = '  -2.40%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E46').Value = '  -1.94%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E47').Value = '  +2.25%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('E48').Value = '  -2.49%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('E49').Value = '  +1.50%  '
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('E50').Value = '  -1.16%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.561.21'
$ws.Range('E51').Value = '  +1.12%  '
